# Update "想去人数" (F column) values on sheet "展览" and sheet "全部类型"
# to reflect freshly re-generated numbers (gh-pages output regenerated).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (first sheet) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value  = 3913
$wsExhibit.Range("F9").Value  = 3185
$wsExhibit.Range("F11").Value = 2335
$wsExhibit.Range("F15").Value = 466
$wsExhibit.Range("F19").Value = 351
$wsExhibit.Range("F22").Value = 674
$wsExhibit.Range("F24").Value = 49
$wsExhibit.Range("F27").Value = 139
$wsExhibit.Range("F28").Value = 159
$wsExhibit.Range("F29").Value = 38
$wsExhibit.Range("F31").Value = 70
$wsExhibit.Range("F32").Value = 4380
$wsExhibit.Range("F33").Value = 4220
$wsExhibit.Range("F35").Value = 139
$wsExhibit.Range("F38").Value = 1152

# --- Sheet "全部类型" (fourth sheet, aggregates all types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F10").Value = 3913
$wsAll.Range("F14").Value = 3185
$wsAll.Range("F17").Value = 2335
$wsAll.Range("F23").Value = 351
$wsAll.Range("F26").Value = 674
$wsAll.Range("F28").Value = 49
$wsAll.Range("F30").Value = 159
$wsAll.Range("F32").Value = 70
$wsAll.Range("F34").Value = 4380
$wsAll.Range("F35").Value = 4220
$wsAll.Range("F38").Value = 1152
